# DOMA-11339: add a "Decommissioning date" column (Q) to the meters import
# example, with a sample value in the first data row (Q2), mirroring the
# look of the existing "Automatic" column (P).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the whole "Automatic" column (header + data cells)
# into the new column so fills/borders/number-format match exactly.
$ws.Range("P1:P11").Copy($ws.Range("Q1:Q11"))

# Header text for the new column.
$ws.Range("Q1").Value = "Decommissioning date"

# Sample decommissioning date for the first meter row; the rest of the new
# column stays empty, same as the "Automatic" column does for this sample.
$ws.Range("Q2").Value = "2022-01-25"

# Column O:P are 23.5 wide; give the new column the same width.
$ws.Range("Q1").ColumnWidth = $ws.Range("O1").ColumnWidth
